$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1891
$ws.Range("I98").Value = 1958.6207
$ws.Range("K98").Value = 1958.6207
$ws.Range("M98").Value = -460.6206999999999
$ws.Range("H106").Value = 1499.5714
$ws.Range("I106").Value = 1499.5714
$ws.Range("K106").Value = 1499.5714
$ws.Range("M106").Value = -868.5714
$ws.Range("H122").Value = 1891
$ws.Range("I122").Value = 1958.6207
$ws.Range("K122").Value = 5875.8621
$ws.Range("M122").Value = -3425.8621
$ws.Range("H125").Value = 8336229.5
$ws.Range("I125").Value = 1808.9166
$ws.Range("J125").Value = 11908124
$ws.Range("K125").Value = 16280.2494
$ws.Range("L125").Value = 107173116
$ws.Range("M125").Value = -13820.2494
$ws.Range("N125").Value = -107178036
$ws.Range("H131").Value = 4774.316
$ws.Range("I131").Value = 1793.1
$ws.Range("K131").Value = 5379.299999999999
$ws.Range("M131").Value = -339.2999999999993
$ws.Range("H137").Value = 57349.53
$ws.Range("I137").Value = 94883.94500000001
$ws.Range("J137").Value = 2491.5386
$ws.Range("K137").Value = 284651.835
$ws.Range("L137").Value = 7474.6158
$ws.Range("M137").Value = -282101.835
$ws.Range("N137").Value = -12574.6158
# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5997641.5
$ws.Range("I45").Value = 11067825
$ws.Range("J45").Value = 5606.636
$ws.Range("K45").Value = 11067825
$ws.Range("L45").Value = 5606.636
$ws.Range("M45").Value = -11067448
$ws.Range("N45").Value = -6360.636
$ws.Range("H61").Value = 4828.8096
$ws.Range("I61").Value = 4983
$ws.Range("K61").Value = 4983
$ws.Range("M61").Value = -4771
$ws.Range("H122").Value = 675527.2
$ws.Range("I122").Value = 2028.9474
$ws.Range("J122").Value = 1741899.4
$ws.Range("K122").Value = 6086.8422
$ws.Range("L122").Value = 5225698.199999999
$ws.Range("M122").Value = -3636.8422
$ws.Range("N122").Value = -5230598.199999999
$ws.Range("H132").Value = 2363.973
$ws.Range("I132").Value = 1536.1
$ws.Range("K132").Value = 4608.299999999999
$ws.Range("M132").Value = -2078.299999999999
$ws.Range("H136").Value = 4828.8096
$ws.Range("I136").Value = 4983
$ws.Range("K136").Value = 14949
$ws.Range("M136").Value = -12399
# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21035.74
$ws.Range("I31").Value = 2609.889
$ws.Range("J31").Value = 30248.666
$ws.Range("K31").Value = 2609.889
$ws.Range("L31").Value = 30248.666
$ws.Range("M31").Value = -2314.889
$ws.Range("N31").Value = -30838.666
$ws.Range("H34").Value = 21035.74
$ws.Range("I34").Value = 2609.889
$ws.Range("J34").Value = 30248.666
$ws.Range("K34").Value = 2609.889
$ws.Range("L34").Value = 30248.666
$ws.Range("M34").Value = -2407.889
$ws.Range("N34").Value = -30652.666
$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 8373.579
$ws.Range("I86").Value = 6600.273
$ws.Range("K86").Value = 6600.273
$ws.Range("M86").Value = -5477.273
$ws.Range("H89").Value = 8373.579
$ws.Range("I89").Value = 6600.273
$ws.Range("K89").Value = 33001.365
$ws.Range("M89").Value = -27385.365
$ws.Range("H132").Value = 81100.875
$ws.Range("I132").Value = 64136.375
$ws.Range("K132").Value = 192409.125
$ws.Range("M132").Value = -189879.125
# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 401500.4
$ws.Range("I114").Value = 53
$ws.Range("J114").Value = 477966.56
$ws.Range("K114").Value = 159
$ws.Range("L114").Value = 1433899.68
$ws.Range("M114").Value = 3095
$ws.Range("N114").Value = -1440407.68
$ws.Range("H117").Value = 519.3
$ws.Range("J117").Value = 406.2857
$ws.Range("L117").Value = 1218.8571
$ws.Range("N117").Value = -8102.8571
# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 19000
$ws.Range("J36").Value = 19000
$ws.Range("L36").Value = 19000
$ws.Range("N36").Value = -19970
$ws.Range("H43").Value = 13458.357
$ws.Range("I43").Value = 5419.5
$ws.Range("J43").Value = 19487.5
$ws.Range("K43").Value = 5419.5
$ws.Range("L43").Value = 19487.5
$ws.Range("M43").Value = -5268.5
$ws.Range("N43").Value = -19789.5
$ws.Range("H62").Value = 674999.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 22000
$ws.Range("J63").Value = 22000
$ws.Range("L63").Value = 22000
$ws.Range("N63").Value = -23372
$ws.Range("H65").Value = 674999.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 22000
$ws.Range("J66").Value = 22000
$ws.Range("L66").Value = 66000
$ws.Range("N66").Value = -72864
$ws.Range("H102").Value = 8435483
$ws.Range("I102").Value = 15874172
$ws.Range("J102").Value = 2649835.5
$ws.Range("K102").Value = 15874172
$ws.Range("L102").Value = 2649835.5
$ws.Range("M102").Value = -15872550
$ws.Range("N102").Value = -2653079.5
$ws.Range("H109").Value = 50859.8
$ws.Range("J109").Value = 50859.8
$ws.Range("L109").Value = 50859.8
$ws.Range("N109").Value = -52939.8
$ws.Range("H122").Value = 688559.9399999999
$ws.Range("I122").Value = 991253.4399999999
$ws.Range("K122").Value = 2973760.32
$ws.Range("M122").Value = -2971310.32
$ws.Range("H126").Value = 4956942
$ws.Range("I126").Value = 3790447.2
$ws.Range("J126").Value = 5956794.5
$ws.Range("K126").Value = 11371341.6
$ws.Range("L126").Value = 17870383.5
$ws.Range("M126").Value = -11368871.6
$ws.Range("N126").Value = -17875323.5
# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7697.5835
$ws.Range("I7").Value = 4395.727
$ws.Range("J7").Value = 10491.462
$ws.Range("K7").Value = 4395.727
$ws.Range("L7").Value = 10491.462
$ws.Range("M7").Value = -4283.727
$ws.Range("N7").Value = -10715.462
$ws.Range("H16").Value = 147.86667
$ws.Range("I16").Value = 112
$ws.Range("K16").Value = 112
$ws.Range("M16").Value = 58
$ws.Range("H122").Value = 5672.9565
$ws.Range("I122").Value = 3318.3
$ws.Range("J122").Value = 7484.231
$ws.Range("K122").Value = 9954.900000000001
$ws.Range("L122").Value = 22452.693
$ws.Range("M122").Value = -7504.900000000001
$ws.Range("N122").Value = -27352.693
$ws.Range("H126").Value = 7697.5835
$ws.Range("I126").Value = 4395.727
$ws.Range("J126").Value = 10491.462
$ws.Range("K126").Value = 13187.181
$ws.Range("L126").Value = 31474.386
$ws.Range("M126").Value = -10717.181
$ws.Range("N126").Value = -36414.386
# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7633.625
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 7815.355
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 7815.355
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -9063.355
$ws.Range("H65").Value = 7633.625
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 7815.355
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 39076.77499999999
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -45316.77499999999
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
$ws.Range("H123").Value = 82186
$ws.Range("J123").Value = 82186
$ws.Range("L123").Value = 82186
$ws.Range("N123").Value = -91986
$ws.Range("H126").Value = 1861.95
$ws.Range("I126").Value = 2047.6
$ws.Range("J126").Value = 1305
$ws.Range("K126").Value = 6142.799999999999
$ws.Range("L126").Value = 3915
$ws.Range("M126").Value = -3672.799999999999
$ws.Range("N126").Value = -8855
$ws.Range("H132").Value = 17136890
$ws.Range("I132").Value = 23259628
$ws.Range("K132").Value = 69778884
$ws.Range("M132").Value = -69776354

Write-Output "Applied all Hyperion_Profits updates"